# Update the BDD step text for cell C3 (shared-string rich text):
# old: When Please ensure you correctly enter "First Name" before moving to the next field.
# new: When you correctly enter "First Name" before moving to the next field.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("C3")

$run1 = 'When you correctly enter'
$run2 = ' "First Name" '
$run3 = 'before moving to the next field.'

$cell.Value = $run1 + $run2 + $run3

$start2 = $run1.Length + 1
$len2 = $run2.Length
$start3 = $start2 + $len2
$len3 = $run3.Length

# Middle run: regular (non-bold) Calibri 11
$chars2 = $cell.Characters($start2, $len2)
$chars2.Font.Bold = $false
$chars2.Font.Name = 'Calibri'
$chars2.Font.Size = 11

# Trailing run: bold Calibri 11 (matches the rest of the cell's base style)
$chars3 = $cell.Characters($start3, $len3)
$chars3.Font.Bold = $true
$chars3.Font.Name = 'Calibri'
$chars3.Font.Size = 11

# The edited cell becomes the active selection.
$ws.Activate()
$cell.Select()
